$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new label strings first, in the order they were originally
# authored, so the shared-strings table grows in the same sequence.
$ws.Range("B17").Value = "total time"
$ws.Range("B18").Value = "night time"
$ws.Range("B19").Value = "day time"
$ws.Range("B21").Value = "nap"
$ws.Range("C16").Value = "walkng total"
$ws.Range("C17").Value = "driving total"
$ws.Range("B22").Value = "walking"
$ws.Range("B23").Value = "driving"

# Row 16
$ws.Range("A16").Value = 0.74513888888888891
$ws.Range("A16").NumberFormat = "h:mm"
$ws.Range("D16").Formula = "=D2+D3+D7+D14"

# Row 17
$ws.Range("A17").Formula = "=A16-A2"
$ws.Range("A17").NumberFormat = "h:mm"
$ws.Range("D17").Formula = "=D5+D6+D10+D11+D13"

# Row 18
$ws.Range("A18").Value = 0.010416666666666666
$ws.Range("A18").NumberFormat = "h:mm"

# Row 19
$ws.Range("A19").Value = 0.52083333333333337
$ws.Range("A19").NumberFormat = "h:mm"

# Row 21
$ws.Range("A21").Value = 0.041666666666666664
$ws.Range("A21").NumberFormat = "h:mm"

# Row 22
$ws.Range("A22").Value = 0.40625
$ws.Range("A22").NumberFormat = "h:mm"

# Row 23
$ws.Range("A23").Value = 0.07291666666666667
$ws.Range("A23").NumberFormat = "h:mm"

# Row 24 (blank row, time-formatted like the rest of column A)
$ws.Range("A24").NumberFormat = "h:mm"

# Update selection to D17
$ws.Range("D17").Select()
